$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCellValue {
    param($cell, $val)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextCellValue $ws.Cells.Item(2, 4) '50.721.75'
Set-TextCellValue $ws.Cells.Item(2, 5) '  -1.05%  '
Set-TextCellValue $ws.Cells.Item(3, 4) '2.919.43'
Set-TextCellValue $ws.Cells.Item(3, 5) '  -1.82%  '
Set-TextCellValue $ws.Cells.Item(4, 4) '0.999'
Set-TextCellValue $ws.Cells.Item(4, 5) '  -0.05%  '
Set-TextCellValue $ws.Cells.Item(5, 4) '374.15'
Set-TextCellValue $ws.Cells.Item(5, 5) '  -1.97%  '
Set-TextCellValue $ws.Cells.Item(6, 4) '99.54'
Set-TextCellValue $ws.Cells.Item(6, 5) '  -2.71%  '
Set-TextCellValue $ws.Cells.Item(7, 4) '0.533'
Set-TextCellValue $ws.Cells.Item(7, 5) '  -1.81%  '
Set-TextCellValue $ws.Cells.Item(8, 5) '  -0.06%  '
Set-TextCellValue $ws.Cells.Item(9, 4) '0.582'
Set-TextCellValue $ws.Cells.Item(9, 5) '  -1.39%  '
Set-TextCellValue $ws.Cells.Item(10, 4) '35.77'
Set-TextCellValue $ws.Cells.Item(10, 5) '  -3.06%  '
Set-TextCellValue $ws.Cells.Item(11, 5) '  -0.70%  '
Set-TextCellValue $ws.Cells.Item(12, 4) '0.0839'
Set-TextCellValue $ws.Cells.Item(12, 5) '  -0.26%  '
Set-TextCellValue $ws.Cells.Item(13, 4) '3.372.94'
Set-TextCellValue $ws.Cells.Item(13, 5) '  -2.07%  '
Set-TextCellValue $ws.Cells.Item(14, 4) '17.92'
Set-TextCellValue $ws.Cells.Item(14, 5) '  -1.34%  '
Set-TextCellValue $ws.Cells.Item(15, 4) '7.51'
Set-TextCellValue $ws.Cells.Item(15, 5) '  +0.09%  '
Set-TextCellValue $ws.Cells.Item(16, 2) 'Uniswap'
Set-TextCellValue $ws.Cells.Item(16, 3) 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCellValue $ws.Cells.Item(16, 4) '11.19'
Set-TextCellValue $ws.Cells.Item(16, 5) '  +52.02%  '
Set-TextCellValue $ws.Cells.Item(17, 2) 'WrappedEther'
Set-TextCellValue $ws.Cells.Item(17, 3) 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCellValue $ws.Cells.Item(17, 4) '2.911.73'
Set-TextCellValue $ws.Cells.Item(17, 5) '  -2.07%  '
Set-TextCellValue $ws.Cells.Item(18, 4) '0.982'
Set-TextCellValue $ws.Cells.Item(18, 5) '  -1.44%  '
Set-TextCellValue $ws.Cells.Item(19, 4) '50.618.45'
Set-TextCellValue $ws.Cells.Item(19, 5) '  -1.14%  '
Set-TextCellValue $ws.Cells.Item(20, 4) '3.04'
Set-TextCellValue $ws.Cells.Item(20, 5) '  -6.49%  '
Set-TextCellValue $ws.Cells.Item(21, 4) '12.27'
Set-TextCellValue $ws.Cells.Item(21, 5) '  -4.08%  '
Set-TextCellValue $ws.Cells.Item(22, 4) '0.0₃0949'
Set-TextCellValue $ws.Cells.Item(22, 5) '  -0.68%  '
Set-TextCellValue $ws.Cells.Item(23, 4) '68.62'
Set-TextCellValue $ws.Cells.Item(23, 5) '  -0.38%  '
Set-TextCellValue $ws.Cells.Item(24, 4) '264.12'
Set-TextCellValue $ws.Cells.Item(24, 5) '  +1.49%  '
Set-TextCellValue $ws.Cells.Item(25, 4) '3.13'
Set-TextCellValue $ws.Cells.Item(25, 5) '  +8.35%  '
Set-TextCellValue $ws.Cells.Item(26, 4) '7.94'
Set-TextCellValue $ws.Cells.Item(26, 5) '  -2.14%  '
Set-TextCellValue $ws.Cells.Item(27, 4) '7.32'
Set-TextCellValue $ws.Cells.Item(27, 5) '  -1.71%  '
Set-TextCellValue $ws.Cells.Item(28, 5) '  +0.04%  '
Set-TextCellValue $ws.Cells.Item(29, 4) '25.38'
Set-TextCellValue $ws.Cells.Item(29, 5) '  -1.88%  '
Set-TextCellValue $ws.Cells.Item(30, 4) '0.162'
Set-TextCellValue $ws.Cells.Item(30, 5) '  -4.02%  '
Set-TextCellValue $ws.Cells.Item(31, 4) '0.109'
Set-TextCellValue $ws.Cells.Item(31, 5) '  -7.66%  '
Set-TextCellValue $ws.Cells.Item(32, 4) '9.93'
Set-TextCellValue $ws.Cells.Item(32, 5) '  +0.94%  '
Set-TextCellValue $ws.Cells.Item(33, 4) '50.71'
Set-TextCellValue $ws.Cells.Item(33, 5) '  -0.23%  '
Set-TextCellValue $ws.Cells.Item(34, 5) '  -0.82%  '
Set-TextCellValue $ws.Cells.Item(35, 4) '32.86'
Set-TextCellValue $ws.Cells.Item(35, 5) '  -4.33%  '
Set-TextCellValue $ws.Cells.Item(36, 4) '0.0435'
Set-TextCellValue $ws.Cells.Item(36, 5) '  -4.01%  '
Set-TextCellValue $ws.Cells.Item(37, 5) '  -0.04%  '
Set-TextCellValue $ws.Cells.Item(38, 4) '3.06'
Set-TextCellValue $ws.Cells.Item(38, 5) '  +2.87%  '
Set-TextCellValue $ws.Cells.Item(39, 4) '0.115'
Set-TextCellValue $ws.Cells.Item(39, 5) '  -0.51%  '
Set-TextCellValue $ws.Cells.Item(40, 4) '16.35'
Set-TextCellValue $ws.Cells.Item(40, 5) '  -3.50%  '
Set-TextCellValue $ws.Cells.Item(41, 4) '1.79'
Set-TextCellValue $ws.Cells.Item(41, 5) '  -1.60%  '
Set-TextCellValue $ws.Cells.Item(42, 4) '2.44'
Set-TextCellValue $ws.Cells.Item(42, 5) '  -4.76%  '
Set-TextCellValue $ws.Cells.Item(43, 4) '119.12'
Set-TextCellValue $ws.Cells.Item(43, 5) '  -3.07%  '
Set-TextCellValue $ws.Cells.Item(44, 4) '20.83'
Set-TextCellValue $ws.Cells.Item(44, 5) '  -2.70%  '
Set-TextCellValue $ws.Cells.Item(45, 5) '  -2.04%  '
Set-TextCellValue $ws.Cells.Item(46, 4) '3.34'
Set-TextCellValue $ws.Cells.Item(46, 5) '  +2.23%  '
Set-TextCellValue $ws.Cells.Item(47, 5) '  -1.23%  '
Set-TextCellValue $ws.Cells.Item(48, 4) '0.266'
Set-TextCellValue $ws.Cells.Item(48, 5) '  -2.54%  '
Set-TextCellValue $ws.Cells.Item(49, 4) '1.977.90'
Set-TextCellValue $ws.Cells.Item(49, 5) '  -2.41%  '
Set-TextCellValue $ws.Cells.Item(50, 4) '0.0323'
Set-TextCellValue $ws.Cells.Item(50, 5) '  -2.34%  '
Set-TextCellValue $ws.Cells.Item(51, 4) '5.16'
Set-TextCellValue $ws.Cells.Item(51, 5) '  +1.26%  '
